# Generate Report for Handoff
# Rename the source-file GUID-based identifiers and bump the handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "d16fb887-9a27-4aeb-b554-7e40fbe96b6e"
$newGuid = "2b92aafd-072d-474c-bfb9-9b19fa30b1e0"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldZhXlf = "$oldGuid.25c74f0096fefb4c55f72b270178a417d4c2a0b5.zh-cn.xlf"
$newZhXlf = "$newGuid.938bf64c4ee846b80da5f81dc77114f694c4c031.zh-cn.xlf"

$oldDeXlf = "$oldGuid.25c74f0096fefb4c55f72b270178a417d4c2a0b5.de-de.xlf"
$newDeXlf = "$newGuid.938bf64c4ee846b80da5f81dc77114f694c4c031.de-de.xlf"

$oldZhTime = "2016-03-07 10:13:14"
$newZhTime = "2016-03-07 10:14:07"

$oldDeTime = "2016-03-07 10:13:28"
$newDeTime = "2016-03-07 10:14:21"

$sheetOverview = $wb.Worksheets.Item("Overview")
$sheetZh = $wb.Worksheets.Item("zh-cn")
$sheetDe = $wb.Worksheets.Item("de-de")

# Overview!A2 - update hyperlink text, keep target address
$h = $sheetOverview.Range("A2").Hyperlinks.Item(1)
$h.TextToDisplay = $newMd
$sheetOverview.Range("A2").Value = $newMd

# zh-cn!A2
$h = $sheetZh.Range("A2").Hyperlinks.Item(1)
$h.TextToDisplay = $newMd
$sheetZh.Range("A2").Value = $newMd

# zh-cn!C2
$h = $sheetZh.Range("C2").Hyperlinks.Item(1)
$h.TextToDisplay = $newZhXlf
$sheetZh.Range("C2").Value = $newZhXlf

# zh-cn!D2 - timestamp string (stored as text)
$sheetZh.Range("D2").Value = $newZhTime

# de-de!A2
$h = $sheetDe.Range("A2").Hyperlinks.Item(1)
$h.TextToDisplay = $newMd
$sheetDe.Range("A2").Value = $newMd

# de-de!C2
$h = $sheetDe.Range("C2").Hyperlinks.Item(1)
$h.TextToDisplay = $newDeXlf
$sheetDe.Range("C2").Value = $newDeXlf

# de-de!D2 - timestamp string (stored as text)
$sheetDe.Range("D2").Value = $newDeTime
